# feat: add 2022-Q1 data
#
# The workbook has sheets: 2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计
# We need to:
#   1. Insert a new "2022-Q1" sheet (with the per-fund holdings table) right
#      before the "总计" (totals) sheet - by repurposing the existing "总计"
#      sheet (renaming + rewriting its contents), so the new sheet keeps the
#      original sheet identity/position.
#   2. Create a brand new "总计" sheet after it, containing the same
#      aggregate table as before, plus a new leading row for 2022-Q1
#      (10 funds held, 3.36 billion yuan market value).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: turn the existing "总计" sheet into the new "2022-Q1" sheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# Use another per-fund sheet (2021-Q1) as the formatting template: it has
# the exact same header row + column-A numbering style we need (style
# index 2: bold, centered, thin-bordered).
$tpl = $wb.Worksheets.Item("2021-Q1")

$tpl.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$tpl.Range("A2").Copy()
$q1.Range("A2:A11").PasteSpecial(-4122)

# Header row
$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

# Data rows (A = running index 0..9, B..G are kept as text like the source
# data, H is numeric)
$q1Rows = @(
    @(0, "000011", "华夏大盘精选混合", "42.34", "92.34", "5.21", "2.2059", 7),
    @(1, "160813", "长盛同盛成长优选灵活配置混合 (LOF)", "12.92", "77.90", "3.01", "0.3889", 7),
    @(2, "519039", "长盛同德主题混合", "12.19", "79.80", "3.08", "0.3755", 7),
    @(3, "630010", "华商价值精选混合", "4.93", "89.49", "3.77", "0.1859", 6),
    @(4, "010155", "长盛核心成长混合A", "3.09", "75.19", "3.28", "0.1014", 7),
    @(5, "630006", "华商产业升级混合", "0.98", "87.95", "3.70", "0.0363", 6),
    @(6, "000057", "中银消费主题混合", "0.80", "86.22", "2.98", "0.0238", 10),
    @(7, "005826", "华夏潜龙精选股票", "0.71", "87.98", "2.99", "0.0212", 9),
    @(8, "005161", "华商上游产业股票", "0.36", "89.02", "3.68", "0.0132", 3),
    @(9, "010156", "长盛核心成长混合C", "0.37", "75.19", "3.28", "0.0121", 7)
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Cells.Item($r,1).Value = $row[0]

    # Fund code, size/position/rank percentages are kept as literal text
    # (not numbers) in the source data - e.g. fund codes like "000011"
    # would lose their leading zeros otherwise. Force text entry via the
    # "@" number format, then strip the format again so the cell ends up
    # with default (no explicit) style, same as the source file.
    $q1.Cells.Item($r,2).NumberFormat = "@"
    $q1.Cells.Item($r,2).Value = $row[1]
    $q1.Cells.Item($r,2).ClearFormats()

    $q1.Cells.Item($r,3).Value = $row[2]

    $q1.Cells.Item($r,4).NumberFormat = "@"
    $q1.Cells.Item($r,4).Value = $row[3]
    $q1.Cells.Item($r,4).ClearFormats()

    $q1.Cells.Item($r,5).NumberFormat = "@"
    $q1.Cells.Item($r,5).Value = $row[4]
    $q1.Cells.Item($r,5).ClearFormats()

    $q1.Cells.Item($r,6).NumberFormat = "@"
    $q1.Cells.Item($r,6).Value = $row[5]
    $q1.Cells.Item($r,6).ClearFormats()

    $q1.Cells.Item($r,7).NumberFormat = "@"
    $q1.Cells.Item($r,7).Value = $row[6]
    $q1.Cells.Item($r,7).ClearFormats()

    $q1.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: create the new "总计" sheet after "2022-Q1"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"
# Touching Outline ensures the sheet gets the same sheetPr/outlinePr/
# pageSetUpPr block that every other sheet in this workbook has.
$dummy = $total.Outline.SummaryBelow

$tpl.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$tpl.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Cells.Item(1,2).Value = "日期"
$total.Cells.Item(1,3).Value = "持有数量(只)"
$total.Cells.Item(1,4).Value = "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 10, 3.36),
    @(1, "2021-Q4", 15, 4.96),
    @(2, "2021-Q3", 4, 0.21),
    @(3, "2021-Q2", 1, 1.99),
    @(4, "2021-Q1", 1, 2.3),
    @(5, "2020-Q4", 10, 2.25)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r,1).Value = $row[0]
    $total.Cells.Item($r,2).Value = $row[1]
    $total.Cells.Item($r,3).Value = $row[2]
    $total.Cells.Item($r,4).Value = $row[3]
    $r = $r + 1
}

Write-Host "done"
